# Daily "cryptos" price-table refresh (coinranking.com scrape), as produced
# by the GitHub Actions bot. Updates Price (col D) / Volume(1h) (col E) for
# most rows, and for a few rows the ranking reshuffled so Coin/Link/Price/
# Volume were all replaced with the entry that's now at that rank.
#
# Column D holds price text like "61.932.06" or "0.0000215" that LOOKS
# numeric to Excel's smart paste. Setting .Value directly would silently
# convert it to a float (losing the literal formatting / trailing zeros),
# so for D we prefix with a leading apostrophe to force text, then reset
# the cell Style back to Normal so we don't leave a stray quote-prefix
# number format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'61.932.06"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +0.27%  '
$ws.Cells.Item(3, 4).Value = "'3.434.75"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +1.18%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = "'408.96"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.83%  '
$ws.Cells.Item(6, 4).Value = "'128.64"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.57%  '
$ws.Cells.Item(7, 5).Value = '  +5.59%  '
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
$ws.Cells.Item(9, 4).Value = "'0.738"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +10.33%  '
$ws.Cells.Item(10, 4).Value = "'0.145"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +19.17%  '
$ws.Cells.Item(11, 4).Value = "'42.65"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.52%  '
$ws.Cells.Item(12, 2).Value = 'ShibaInu'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(12, 4).Value = "'0.0000215"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +67.30%  '
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).Value = "'0.141"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -0.39%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = "'3.966.24"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +1.25%  '
$ws.Cells.Item(15, 2).Value = 'Chainlink'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(15, 4).Value = "'21.27"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +7.29%  '
$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(16, 4).Value = "'8.90"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +5.54%  '
$ws.Cells.Item(17, 4).Value = "'3.420.40"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.31%  '
$ws.Cells.Item(18, 4).Value = "'12.40"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +12.69%  '
$ws.Cells.Item(19, 5).Value = '  +4.99%  '
$ws.Cells.Item(20, 4).Value = "'61.888.36"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.35%  '
$ws.Cells.Item(21, 4).Value = "'401.18"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +27.45%  '
$ws.Cells.Item(22, 4).Value = "'89.97"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +5.41%  '
$ws.Cells.Item(23, 5).Value = '  -0.09%  '
$ws.Cells.Item(24, 4).Value = "'13.40"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +4.97%  '
$ws.Cells.Item(25, 4).Value = "'3.21"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +2.69%  '
$ws.Cells.Item(26, 4).Value = "'32.98"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +11.60%  '
$ws.Cells.Item(27, 4).Value = "'8.67"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +4.27%  '
$ws.Cells.Item(28, 4).Value = "'4.81"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +0.36%  '
$ws.Cells.Item(29, 4).Value = "'7.60"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.92%  '
$ws.Cells.Item(30, 4).Value = "'2.72"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +1.71%  '
$ws.Cells.Item(31, 5).Value = '  +2.08%  '
$ws.Cells.Item(32, 5).Value = '  +0.46%  '
$ws.Cells.Item(33, 4).Value = "'11.89"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +4.59%  '
$ws.Cells.Item(34, 4).Value = "'43.63"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +4.98%  '
$ws.Cells.Item(35, 5).Value = '  -0.17%  '
$ws.Cells.Item(36, 4).Value = "'0.0506"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +5.44%  '
$ws.Cells.Item(37, 4).Value = "'53.94"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +4.16%  '
$ws.Cells.Item(38, 5).Value = '  +0.08%  '
$ws.Cells.Item(39, 5).Value = '  -1.07%  '
$ws.Cells.Item(40, 2).Value = 'Stellar'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(40, 4).Value = "'0.132"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +6.22%  '
$ws.Cells.Item(41, 2).Value = 'Stacks'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(41, 4).Value = "'2.91"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -1.24%  '
$ws.Cells.Item(42, 5).Value = '  +6.64%  '
$ws.Cells.Item(43, 4).Value = "'142.03"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +1.97%  '
$ws.Cells.Item(44, 5).Value = '  -0.25%  '
$ws.Cells.Item(45, 4).Value = "'4.05"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +1.56%  '
$ws.Cells.Item(46, 4).Value = "'2.41"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +8.40%  '
$ws.Cells.Item(47, 4).Value = "'16.68"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -0.02%  '
$ws.Cells.Item(48, 4).Value = "'21.73"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.47%  '
$ws.Cells.Item(49, 4).Value = "'2.122.09"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.13%  '
$ws.Cells.Item(50, 5).Value = '  +14.91%  '
$ws.Cells.Item(51, 5).Value = '  +7.95%  '
